$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply all cell text updates from the cryptos list refresh.
# NumberFormat is forced to text ("@") before assignment and the cell
# style is reset back to "Normal" afterward so that values which look
# numeric (e.g. "581.10") are still stored as text, matching the source data.
function Set-TextValue($ws, $ref, $val) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws "D2" "63.739.80"
Set-TextValue $ws "E2" "  -0.85%  "
Set-TextValue $ws "D3" "3.425.43"
Set-TextValue $ws "E3" "  -2.00%  "
Set-TextValue $ws "E4" "  -0.01%  "
Set-TextValue $ws "D5" "581.10"
Set-TextValue $ws "D6" "129.81"
Set-TextValue $ws "E6" "  -3.24%  "
Set-TextValue $ws "E7" "  -0.05%  "
Set-TextValue $ws "D8" "0.481"
Set-TextValue $ws "E8" "  -1.39%  "
Set-TextValue $ws "E9" "  +4.21%  "
Set-TextValue $ws "E10" "  +0.58%  "
Set-TextValue $ws "E11" "  -0.55%  "
Set-TextValue $ws "D12" "4.004.13"
Set-TextValue $ws "E12" "  -2.15%  "
Set-TextValue $ws "E13" "  -0.26%  "
Set-TextValue $ws "E14" "  -1.75%  "
Set-TextValue $ws "D15" "3.424.70"
Set-TextValue $ws "E15" "  -2.05%  "
Set-TextValue $ws "D16" "63.763.44"
Set-TextValue $ws "E16" "  -0.94%  "
Set-TextValue $ws "D17" "24.99"
Set-TextValue $ws "E17" "  -2.66%  "
Set-TextValue $ws "D18" "9.88"
Set-TextValue $ws "E18" "  +0.26%  "
Set-TextValue $ws "D19" "5.68"
Set-TextValue $ws "E19" "  -1.15%  "
Set-TextValue $ws "E20" "  -1.53%  "
Set-TextValue $ws "D21" "384.72"
Set-TextValue $ws "E21" "  -2.20%  "
Set-TextValue $ws "D22" "0.565"
Set-TextValue $ws "E22" "  -1.14%  "
Set-TextValue $ws "D23" "3.562.33"
Set-TextValue $ws "E23" "  -2.05%  "
Set-TextValue $ws "D24" "73.78"
Set-TextValue $ws "E24" "  -1.18%  "
Set-TextValue $ws "E25" "  +0.18%  "
Set-TextValue $ws "D26" "0.0000111"
Set-TextValue $ws "E26" "  -4.49%  "
Set-TextValue $ws "D27" "0.998"
Set-TextValue $ws "E27" "  -0.12%  "
Set-TextValue $ws "D28" "2.20"
Set-TextValue $ws "E28" "  -1.53%  "
Set-TextValue $ws "E29" "  -4.14%  "
Set-TextValue $ws "B30" "InternetComputer(DFINITY)"
Set-TextValue $ws "C30" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws "D30" "7.97"
Set-TextValue $ws "E30" "  -3.33%  "
Set-TextValue $ws "B31" "Kaspa"
Set-TextValue $ws "C31" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws "D31" "0.154"
Set-TextValue $ws "E31" "  +1.46%  "
Set-TextValue $ws "E32" "  -3.59%  "
Set-TextValue $ws "D33" "3.453.61"
Set-TextValue $ws "E33" "  -1.88%  "
Set-TextValue $ws "E34" "  -0.06%  "
Set-TextValue $ws "D35" "22.90"
Set-TextValue $ws "E35" "  -2.30%  "
Set-TextValue $ws "D36" "5.18"
Set-TextValue $ws "E36" "  +0.98%  "
Set-TextValue $ws "D37" "6.78"
Set-TextValue $ws "E37" "  -1.50%  "
Set-TextValue $ws "D38" "163.88"
Set-TextValue $ws "E38" "  -2.09%  "
Set-TextValue $ws "D39" "1.51"
Set-TextValue $ws "E39" "  -2.53%  "
Set-TextValue $ws "D40" "0.0775"
Set-TextValue $ws "E40" "  -0.57%  "
Set-TextValue $ws "D41" "0.790"
Set-TextValue $ws "E41" "  -2.16%  "
Set-TextValue $ws "E42" "  -0.08%  "
Set-TextValue $ws "D43" "41.43"
Set-TextValue $ws "E43" "  -1.08%  "
Set-TextValue $ws "D44" "4.35"
Set-TextValue $ws "E44" "  -1.00%  "
Set-TextValue $ws "E45" "  -2.26%  "
Set-TextValue $ws "E46" "  -7.66%  "
Set-TextValue $ws "E47" "  -3.98%  "
Set-TextValue $ws "D49" "0.900"
Set-TextValue $ws "E49" "  +0.79%  "
Set-TextValue $ws "D50" "2.294.47"
Set-TextValue $ws "E50" "  -6.99%  "
Set-TextValue $ws "E51" "  -2.09%  "
